# Apply updated RandomForest imputation values in column E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E4"  = 16.50069999999999
    "E6"  = 16.4525
    "E7"  = 15.75690000000001
    "E8"  = 16.7053
    "E16" = 16.19510000000001
    "E20" = 15.8219
    "E21" = 16.7891
    "E28" = 16.5269
    "E29" = 17.13820000000001
    "E30" = 15.1561
    "E32" = 17.01099999999998
    "E40" = 16.98539999999999
    "E46" = 17.01979999999999
    "E51" = 17.27610000000001
    "E52" = 17.11820000000001
    "E57" = 16.52580000000001
    "E59" = 15.89310000000001
    "E62" = 16.6693
    "E66" = 17.1634
    "E73" = 17.241
    "E74" = 16.72919999999998
    "E77" = 17.83430000000001
    "E92" = 18.50270000000002
    "E100" = 16.3545
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
